# Apply the two content edits to
# "Day 2/Hands On Assignments/2. Working with Collections/Hands On Exercises - Map.docx"

$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# Fix the accidental double period at the end of the first "NOTE" line
# (Assignment 02 / TreeMap section).
$d.Content.Find.Execute(
    "NOTE: You can test the methods using a main method..",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "NOTE: You can test the methods using a main method.",
    2) | Out-Null

# --- Edit 2 -----------------------------------------------------------
# Assignment 06 (HashTable) intro line was missing the assignment number:
# "Implement the assignment  using HashTable" (note the double space)
# becomes "Implement the assignment 1 using HashTable".
$r = $d.Content
$r.Find.Execute(
    "Implement the assignment  using HashTable",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null

# Position right after "Implement the assignment " (25 chars in) and type "1".
$ins = $r.Duplicate
$ins.Collapse(1)
$ins.MoveStart(1, 25) | Out-Null
$ins.Collapse(1)
$ins.InsertAfter("1")

# Word re-anchors its hidden "_GoBack" bookmark at the point of the last
# edit, so move it to sit right after the newly typed "1".
$bm = $ins.Duplicate
$bm.Collapse(1)
$bm.MoveStart(1, 1) | Out-Null
$bm.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null
